# Fixed update to excel issue
# 1. Rename headers on existing sheets (Weekly Quantity, Monthly Trend)
# 2. Add a new "PO Forecast" sheet with forecast data

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# Update the "Requested quantity" headers to their new names
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" worksheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "PO Forecast"

# Copy the header-row formatting (bold, centered, bordered) from the
# Weekly Quantity sheet so the new header row looks consistent
$wsWeekly.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# Copy the date-formatted style used in column A down the forecast rows
$wsWeekly.Range("A2").Copy()
$ws3.Range("A2:A31").PasteSpecial(-4122)

# Populate the forecast data (ds, PO_Forecast, yhat_lower, yhat_upper)
$ws3.Cells.Item(2, 1).Value = 44955.99999999999
$ws3.Cells.Item(2, 2).Value = 60
$ws3.Cells.Item(2, 3).Value = -49.22020192401121
$ws3.Cells.Item(2, 4).Value = 182.3989882141961
$ws3.Cells.Item(3, 1).Value = 44983.99999999999
$ws3.Cells.Item(3, 2).Value = 69
$ws3.Cells.Item(3, 3).Value = -42.26879408932941
$ws3.Cells.Item(3, 4).Value = 182.9313087250131
$ws3.Cells.Item(4, 1).Value = 45039.99999999999
$ws3.Cells.Item(4, 2).Value = 87
$ws3.Cells.Item(4, 3).Value = -32.09473359230542
$ws3.Cells.Item(4, 4).Value = 203.2629513878829
$ws3.Cells.Item(5, 1).Value = 45046.99999999999
$ws3.Cells.Item(5, 2).Value = 89
$ws3.Cells.Item(5, 3).Value = -22.3395571409629
$ws3.Cells.Item(5, 4).Value = 206.8440167880495
$ws3.Cells.Item(6, 1).Value = 45053.99999999999
$ws3.Cells.Item(6, 2).Value = 91
$ws3.Cells.Item(6, 3).Value = -29.01641171759038
$ws3.Cells.Item(6, 4).Value = 208.4707269691056
$ws3.Cells.Item(7, 1).Value = 45060.99999999999
$ws3.Cells.Item(7, 2).Value = 93
$ws3.Cells.Item(7, 3).Value = -20.3280500619218
$ws3.Cells.Item(7, 4).Value = 205.1946712407021
$ws3.Cells.Item(8, 1).Value = 45067.99999999999
$ws3.Cells.Item(8, 2).Value = 96
$ws3.Cells.Item(8, 3).Value = -21.8493264558648
$ws3.Cells.Item(8, 4).Value = 214.0838397922211
$ws3.Cells.Item(9, 1).Value = 45074.99999999999
$ws3.Cells.Item(9, 2).Value = 98
$ws3.Cells.Item(9, 3).Value = -17.14428278668671
$ws3.Cells.Item(9, 4).Value = 218.2961612963763
$ws3.Cells.Item(10, 1).Value = 45081.99999999999
$ws3.Cells.Item(10, 2).Value = 100
$ws3.Cells.Item(10, 3).Value = -17.1981175666288
$ws3.Cells.Item(10, 4).Value = 218.7236708130193
$ws3.Cells.Item(11, 1).Value = 45088.99999999999
$ws3.Cells.Item(11, 2).Value = 102
$ws3.Cells.Item(11, 3).Value = -4.015692686294443
$ws3.Cells.Item(11, 4).Value = 212.8672626827249
$ws3.Cells.Item(12, 1).Value = 45102.99999999999
$ws3.Cells.Item(12, 2).Value = 107
$ws3.Cells.Item(12, 3).Value = -5.038645663146046
$ws3.Cells.Item(12, 4).Value = 226.3696306682118
$ws3.Cells.Item(13, 1).Value = 45109.99999999999
$ws3.Cells.Item(13, 2).Value = 109
$ws3.Cells.Item(13, 3).Value = 3.694900871599079
$ws3.Cells.Item(13, 4).Value = 219.8060509269497
$ws3.Cells.Item(14, 1).Value = 45116.99999999999
$ws3.Cells.Item(14, 2).Value = 111
$ws3.Cells.Item(14, 3).Value = -4.461209944971922
$ws3.Cells.Item(14, 4).Value = 232.9423296548443
$ws3.Cells.Item(15, 1).Value = 45123.99999999999
$ws3.Cells.Item(15, 2).Value = 113
$ws3.Cells.Item(15, 3).Value = 11.05491669316999
$ws3.Cells.Item(15, 4).Value = 232.1510194981597
$ws3.Cells.Item(16, 1).Value = 45130.99999999999
$ws3.Cells.Item(16, 2).Value = 116
$ws3.Cells.Item(16, 3).Value = -2.866600220316528
$ws3.Cells.Item(16, 4).Value = 227.6102392223659
$ws3.Cells.Item(17, 1).Value = 45137.99999999999
$ws3.Cells.Item(17, 2).Value = 118
$ws3.Cells.Item(17, 3).Value = 5.717740374423427
$ws3.Cells.Item(17, 4).Value = 235.6978991832153
$ws3.Cells.Item(18, 1).Value = 45151.99999999999
$ws3.Cells.Item(18, 2).Value = 122
$ws3.Cells.Item(18, 3).Value = -3.9527125234518
$ws3.Cells.Item(18, 4).Value = 243.3987964573794
$ws3.Cells.Item(19, 1).Value = 45158.99999999999
$ws3.Cells.Item(19, 2).Value = 124
$ws3.Cells.Item(19, 3).Value = 8.456903623963219
$ws3.Cells.Item(19, 4).Value = 241.8652210678987
$ws3.Cells.Item(20, 1).Value = 45165.99999999999
$ws3.Cells.Item(20, 2).Value = 127
$ws3.Cells.Item(20, 3).Value = 11.1317532853422
$ws3.Cells.Item(20, 4).Value = 236.6924273813061
$ws3.Cells.Item(21, 1).Value = 45172.99999999999
$ws3.Cells.Item(21, 2).Value = 129
$ws3.Cells.Item(21, 3).Value = 19.00140511592284
$ws3.Cells.Item(21, 4).Value = 252.1574521163286
$ws3.Cells.Item(22, 1).Value = 45179.99999999999
$ws3.Cells.Item(22, 2).Value = 131
$ws3.Cells.Item(22, 3).Value = 23.81287922497725
$ws3.Cells.Item(22, 4).Value = 255.704728175635
$ws3.Cells.Item(23, 1).Value = 45193.99999999999
$ws3.Cells.Item(23, 2).Value = 136
$ws3.Cells.Item(23, 3).Value = 25.44757615294709
$ws3.Cells.Item(23, 4).Value = 254.6801257883764
$ws3.Cells.Item(24, 1).Value = 45200.99999999999
$ws3.Cells.Item(24, 2).Value = 138
$ws3.Cells.Item(24, 3).Value = 12.66890904049576
$ws3.Cells.Item(24, 4).Value = 244.6122799852899
$ws3.Cells.Item(25, 1).Value = 45207.99999999999
$ws3.Cells.Item(25, 2).Value = 140
$ws3.Cells.Item(25, 3).Value = 30.99888066513418
$ws3.Cells.Item(25, 4).Value = 256.9757762652324
$ws3.Cells.Item(26, 1).Value = 45214.99999999999
$ws3.Cells.Item(26, 2).Value = 142
$ws3.Cells.Item(26, 3).Value = 22.38080956464855
$ws3.Cells.Item(26, 4).Value = 261.3322169180728
$ws3.Cells.Item(27, 1).Value = 45221.99999999999
$ws3.Cells.Item(27, 2).Value = 144
$ws3.Cells.Item(27, 3).Value = 29.39910890823109
$ws3.Cells.Item(27, 4).Value = 258.18437119269
$ws3.Cells.Item(28, 1).Value = 45228.99999999999
$ws3.Cells.Item(28, 2).Value = 147
$ws3.Cells.Item(28, 3).Value = 30.32978074991449
$ws3.Cells.Item(28, 4).Value = 263.6084442392992
$ws3.Cells.Item(29, 1).Value = 45235.99999999999
$ws3.Cells.Item(29, 2).Value = 149
$ws3.Cells.Item(29, 3).Value = 34.86100017740583
$ws3.Cells.Item(29, 4).Value = 262.31811152266
$ws3.Cells.Item(30, 1).Value = 45242.99999999999
$ws3.Cells.Item(30, 2).Value = 151
$ws3.Cells.Item(30, 3).Value = 36.38190495625549
$ws3.Cells.Item(30, 4).Value = 263.3742130259569
$ws3.Cells.Item(31, 1).Value = 45249.99999999999
$ws3.Cells.Item(31, 2).Value = 153
$ws3.Cells.Item(31, 3).Value = 37.26006441422233
$ws3.Cells.Item(31, 4).Value = 271.3080613512799

$ws3.Range("A1").Select()
